$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9810446500778198
$ws.Range("B1").Value = 1.95371150970459
$ws.Range("C1").Value = 5.22498893737793
$ws.Range("D1").Value = 1.384364247322083
$ws.Range("E1").Value = 0.8146126270294189
